$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 11 (keep header row 1 and first data row 2)
$ws.Rows("3:11").Delete()

# Update the timestamp in the remaining data row
$ws.Range("C2").Value = "25/07/2025 14:58:08"
